# Efnb2-Epha3.xlsx: refresh with newly computed TPM-based NATMI scores.
# The underlying ligand/receptor TPM values changed, which both updates the
# numeric columns and removes the "ECs" target-cluster rows (ECs is no longer
# a valid receptor-expressing cluster for Epha3 under the new TPM cutoffs), so
# the data block shrinks from 9 rows (3 senders x 3 targets) to 6 rows
# (3 senders x 2 targets: FAPs, MuSCs). Dimension goes from A1:T10 to A1:T7.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the last 3 rows (old rows 8-10, the MuSCs-sender block that included the
# now-removed ECs target-cluster entries alongside FAPs/MuSCs)
$ws.Range("A8:T10").EntireRow.Delete() | Out-Null

# Row 2: ECs -> FAPs
$ws.Range("A2").Value = "ECs"
$ws.Range("D2").Value = "FAPs"
$ws.Range("G2").Value = 36.94436433333333
$ws.Range("H2").Value = 110.833093
$ws.Range("I2").Value = 0.8328964975864823
$ws.Range("J2").Value = 0.8328964975864824
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 26.097779
$ws.Range("N2").Value = 78.29333700000001
$ws.Range("O2").Value = 0.9922055808976035
$ws.Range("P2").Value = 0.9922055808976036
$ws.Range("Q2").Value = 964.1658556668158
$ws.Range("R2").Value = 8677.492701001342
$ws.Range("S2").Value = 0.8264045532153751
$ws.Range("T2").Value = 0.8264045532153753

# Row 3: ECs -> MuSCs
$ws.Range("A3").Value = "ECs"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("G3").Value = 36.94436433333333
$ws.Range("H3").Value = 110.833093
$ws.Range("I3").Value = 0.8328964975864823
$ws.Range("J3").Value = 0.8328964975864824
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.205015
$ws.Range("N3").Value = 0.6150450000000001
$ws.Range("O3").Value = 0.007794419102396499
$ws.Range("P3").Value = 0.007794419102396499
$ws.Range("Q3").Value = 7.574148853798334
$ws.Range("R3").Value = 68.16733968418501
$ws.Range("S3").Value = 0.006491944371107217
$ws.Range("T3").Value = 0.006491944371107218

# Row 4: FAPs -> FAPs
$ws.Range("A4").Value = "FAPs"
$ws.Range("D4").Value = "FAPs"
$ws.Range("G4").Value = 3.374819
$ws.Range("H4").Value = 10.124457
$ws.Range("I4").Value = 0.07608399754092349
$ws.Range("J4").Value = 0.07608399754092349
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 26.097779
$ws.Range("N4").Value = 78.29333700000001
$ws.Range("O4").Value = 0.9922055808976035
$ws.Range("P4").Value = 0.9922055808976036
$ws.Range("Q4").Value = 88.075280427001
$ws.Range("R4").Value = 792.677523843009
$ws.Range("S4").Value = 0.07549096697710382
$ws.Range("T4").Value = 0.07549096697710383

# Row 5: FAPs -> MuSCs
$ws.Range("A5").Value = "FAPs"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("G5").Value = 3.374819
$ws.Range("H5").Value = 10.124457
$ws.Range("I5").Value = 0.07608399754092349
$ws.Range("J5").Value = 0.07608399754092349
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.205015
$ws.Range("N5").Value = 0.6150450000000001
$ws.Range("O5").Value = 0.007794419102396499
$ws.Range("P5").Value = 0.007794419102396499
$ws.Range("Q5").Value = 0.6918885172850001
$ws.Range("R5").Value = 6.226996655565
$ws.Range("S5").Value = 0.0005930305638196622
$ws.Range("T5").Value = 0.0005930305638196622

# Row 6: MuSCs -> FAPs
$ws.Range("A6").Value = "MuSCs"
$ws.Range("D6").Value = "FAPs"
$ws.Range("G6").Value = 4.037305666666668
$ws.Range("H6").Value = 12.111917
$ws.Range("I6").Value = 0.09101950487259411
$ws.Range("J6").Value = 0.09101950487259411
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 26.097779
$ws.Range("N6").Value = 78.29333700000001
$ws.Range("O6").Value = 0.9922055808976035
$ws.Range("P6").Value = 0.9922055808976036
$ws.Range("Q6").Value = 105.3647110441144
$ws.Range("R6").Value = 948.2823993970293
$ws.Range("S6").Value = 0.09031006070512448
$ws.Range("T6").Value = 0.0903100607051245

# Row 7: MuSCs -> MuSCs
$ws.Range("A7").Value = "MuSCs"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("G7").Value = 4.037305666666668
$ws.Range("H7").Value = 12.111917
$ws.Range("I7").Value = 0.09101950487259411
$ws.Range("J7").Value = 0.09101950487259411
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.205015
$ws.Range("N7").Value = 0.6150450000000001
$ws.Range("O7").Value = 0.007794419102396499
$ws.Range("P7").Value = 0.007794419102396499
$ws.Range("Q7").Value = 0.827708221251667
$ws.Range("R7").Value = 7.449373991265002
$ws.Range("S7").Value = 0.0007094441674696187
$ws.Range("T7").Value = 0.0007094441674696187

